$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 & 8: normalize B/C formatting to match the pattern already used
#     by rows 9-12 (pink highlight block), and fill in the missing
#     "Gender formatted" value ("other") in column G ---

$ws.Range("B9").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

$ws.Range("C10").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null

$ws.Range("G7").Value = "other"
$ws.Range("G8").Value = "other"

$excel.CutCopyMode = 0

# --- Row 14: this record actually belongs with the plain (unhighlighted)
#     rows above/below it (13 and 15), so re-format it to match row 15 and
#     drop the stray Gender cell / fill in the formatted phone + gender ---

$ws.Range("A15").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null

$ws.Range("B15").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

$ws.Range("D15").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null

$ws.Range("E15").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("C14").Clear()

$ws.Range("F14").ClearFormats()
$ws.Range("G14").ClearFormats()
$ws.Range("G14").Value = "other"
